$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The instruction packet moved from 13 modeled devices / breakers down to
# 6 (v1.1: 1 voltage bus / 6xCB). Remove the now-unneeded device rows
# (old rows 8-14, i.e. devices #7-#13) - this also shifts the trailing
# blank template rows up so the sheet ends at row 12 instead of row 19.
$ws.Rows("8:14").Delete()

# Reposition the floating notes/legend textbox so it still sits just
# below the (now shorter) data table, anchored under the new row 10.
$shp = $ws.Shapes.Item(1)
$shp.Top = $ws.Range("A10").Top() + 1.5

# Restore the cursor/selection position left by the edit.
$ws.Range("B6").Select()
